# Updates the crypto price/volume table with freshly scraped values.
# Note: many "Price" (column D) values look like plain numbers (e.g. "314.04"),
# but the sheet stores them as text. Assigning such a string to a Range.Value
# would normally make Excel auto-convert it to a floating point number (losing
# the exact text, e.g. "314.04" -> 314.04000000000002). To keep the cell as
# text - matching the original data - we temporarily force the cell's
# NumberFormat to "@" (text) before assigning the value, then restore the
# cell style back to "Normal" so no extra formatting lingers on the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '27.301.20'
$ws.Range("E2").Value = '  +1.19%  '

# Row 3
$ws.Range("D3").Value = '1.852.39'
$ws.Range("E3").Value = '  +1.38%  '

# Row 4
$ws.Range("E4").Value = '  -0.82%  '

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '314.04'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +0.99%  '

# Row 6
$ws.Range("E6").Value = '  -0.63%  '

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4608'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.97%  '

# Row 8
$ws.Range("E8").Value = '  +0.42%  '

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07301'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.68%  '

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.8850'
$ws.Range("D10").Style = "Normal"

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '19.92'
$ws.Range("D11").Style = "Normal"

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07788'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  -0.68%  '

# Row 13
$ws.Range("D13").Value = '1.860.84'
$ws.Range("E13").Value = '  -3.26%  '

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.372'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.68%  '

# Row 15
$ws.Range("E15").Value = '  -0.19%  '

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '91.64'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -0.08%  '

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -0.83%  '

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000008968'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.48%  '

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.000'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.64%  '

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '14.76'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.08%  '

# Row 21
$ws.Range("D21").Value = '27.325.10'
$ws.Range("E21").Value = '  +2.22%  '

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.122'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  -0.07%  '

# Row 23
$ws.Range("E23").Value = '  -0.38%  '

# Row 24
$ws.Range("D24").Value = '2.058.57'
$ws.Range("E24").Value = '  -0.62%  '

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '1.917'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +4.25%  '

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '151.63'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.42%  '

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '18.33'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.18%  '

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '2.058'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -1.37%  '

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '116.17'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +0.91%  '

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '5.092'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +0.59%  '

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.08836'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  -0.44%  '

# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '3.130'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.70%  '

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.7728'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +6.61%  '

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.169'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +3.26%  '

# Row 35
$ws.Range("E35").Value = '  +1.53%  '

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '2.661'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +7.19%  '

# Row 37
$ws.Range("B37").Value = 'TrustWalletToken'
$ws.Range("C37").Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '1.078'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +0.64%  '

# Row 38
$ws.Range("B38").Value = 'VeChain'
$ws.Range("C38").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.01959'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +0.53%  '

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.05234'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.33%  '

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '2.962'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +1.27%  '

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '6.993'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -1.70%  '

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.5141'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.73%  '

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.1632'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.33%  '

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '8.414'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +2.58%  '

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.4811'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.10%  '

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '10.30'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.57%  '

# Row 47
$ws.Range("E47").Value = '  -0.71%  '

# Row 48
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '102.59'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +0.46%  '

# Row 49
$ws.Range("B49").Value = 'NEARProtocol'
$ws.Range("C49").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.647'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +1.47%  '

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.06218'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.05%  '

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '65.54'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +2.48%  '
